$d = $word.ActiveDocument

# ------------------------------------------------------------------
# The "Meta description" paragraph (2nd paragraph of the doc) has the
# exact run pattern we need for the new heading paragraph we are about
# to create near the end of the document: an empty leading run followed
# by a bold run. Grab that run's FormattedText *before* touching
# anything else in the document (FormattedText is a live range, so it
# must be captured and consumed before any other edit shifts things).
# ------------------------------------------------------------------
$metaPara = $d.Paragraphs(2)
$metaStart = $metaPara.Range.Start
$labelLen = "Meta description".Length
$labelRange = $d.Range($metaStart, $metaStart + $labelLen)
$boldRunFormat = $labelRange.FormattedText

# ------------------------------------------------------------------
# Insert a new empty paragraph right before the final paragraph (the
# one holding the italic image-prompt text), then drop the captured
# "<empty-run><bold-run>" pattern into it.
# ------------------------------------------------------------------
$count = $d.Paragraphs.Count
$beforeLast = $d.Paragraphs($count - 1)
$beforeLast.Range.InsertParagraphAfter()

$newCount = $d.Paragraphs.Count
$newParaIndex = $newCount - 1
$newPara = $d.Paragraphs($newParaIndex)
$newPara.Style = "Normal"
$newParaStart = $newPara.Range.Start
$newParaDest = $d.Range($newParaStart, $newParaStart)
$newParaDest.FormattedText = $boldRunFormat

# ------------------------------------------------------------------
# Now remove the original "Meta description" paragraph entirely.
# ------------------------------------------------------------------
$metaPara2 = $d.Paragraphs(2)
$metaPara2.Range.Delete()

# ------------------------------------------------------------------
# Swap the placeholder "Meta description" text (now living in the new
# bold paragraph) for the real heading copy.
# ------------------------------------------------------------------
$d.Content.Find.Execute("Meta description", $true, $false, $false, $false, $false, $true, 1, $false, "Play Columbus slot game for free - Review of Novomatic's historical game", 2)

# ------------------------------------------------------------------
# Replace the closing italic paragraph's text (the old image-generation
# prompt) with the new meta-description copy, keeping its italic run
# formatting untouched.
# ------------------------------------------------------------------
$oldClosing = 'Create a feature image for the Novomatic game "Columbus" with a cartoon-style happy Maya warrior wearing glasses. The image can include elements such as the three caravels, Christopher Columbus, and other symbols from the game, as well as a background that fits the historical theme of the discovery of America. Use bright colors and a fun, engaging design to capture the attention of players and convey the excitement of the game.'
$newClosing = 'Enjoy the historical theme of Columbus and try to win big on the 5-reel, 9-payline slot machine by Novomatic. Play for free.'
$d.Content.Find.Execute($oldClosing, $true, $false, $false, $false, $false, $true, 1, $false, $newClosing, 2)
